# Disable/remove the "Display" (CLK / DIO / Display) sensor block on Plan1.
# The three label cells that described the 4-digit 7-segment display wiring
# (CLK in J11, Display in K11, DIO in J12) are cleared out, while the cell
# formatting (borders/alignment) of the block is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

$ws.Range("J11").Value = ""
$ws.Range("K11").Value = ""
$ws.Range("J12").Value = ""

# Move/save the current selection as it was left by the author after editing.
$ws.Range("N13").Select()
